$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -8.071
$ws.Range("C7").Value = -13.435
$ws.Range("E7").Value = 15.978
$ws.Range("B9").Value = 5.361
$ws.Range("E10").Value = 16.565
$ws.Range("C12").Value = -11.141
$ws.Range("B13").Value = 5.057
$ws.Range("E13").Value = 16.639
$ws.Range("C14").Value = -13.321
$ws.Range("D15").Value = -8.354000000000001
$ws.Range("B16").Value = 5.244999999999999
$ws.Range("E16").Value = 16.675
$ws.Range("B18").Value = 4.795
$ws.Range("C19").Value = -12.028
$ws.Range("B20").Value = 6.775999999999999
$ws.Range("E20").Value = 15.951
$ws.Range("E24").Value = 16.763
$ws.Range("B26").Value = 6.118
$ws.Range("C26").Value = -12.795
$ws.Range("B27").Value = 5.956
$ws.Range("C27").Value = -13.768
$ws.Range("D28").Value = -7.904000000000001
$ws.Range("B29").Value = 5.36
$ws.Range("C29").Value = -11.281
$ws.Range("E32").Value = 16.394
$ws.Range("D33").Value = -7.342000000000001
$ws.Range("B35").Value = 9.239000000000001
$ws.Range("D35").Value = -7.826000000000001
$ws.Range("B36").Value = 8.279
$ws.Range("C37").Value = -13.038
$ws.Range("C38").Value = -14.03
$ws.Range("D38").Value = -8.324999999999999
$ws.Range("E39").Value = 16.532
$ws.Range("D43").Value = -7.776999999999999
$ws.Range("D44").Value = -7.476999999999999
$ws.Range("B45").Value = 6.018
$ws.Range("D45").Value = -7.500999999999999
$ws.Range("C47").Value = -11.881
$ws.Range("D47").Value = -7.132
$ws.Range("E47").Value = 17.259
$ws.Range("E48").Value = 17.15
$ws.Range("C51").Value = -13.014
$ws.Range("D51").Value = -7.57
$ws.Range("C52").Value = -11.99
$ws.Range("E52").Value = 17.191
$ws.Range("D54").Value = -8.084
$ws.Range("B55").Value = 5.757000000000001
$ws.Range("C55").Value = -13.552
$ws.Range("E56").Value = 16.448
$ws.Range("B57").Value = 5.423999999999999
$ws.Range("D57").Value = -8.006
$ws.Range("D62").Value = -7.867
$ws.Range("D63").Value = -7.641999999999999
$ws.Range("D67").Value = -6.986
$ws.Range("B69").Value = 5.293000000000001
$ws.Range("C69").Value = -10.921
$ws.Range("C70").Value = -13.014
$ws.Range("D70").Value = -7.863
$ws.Range("B76").Value = 5.548
$ws.Range("C76").Value = -12.395
$ws.Range("B78").Value = 7.523999999999999
$ws.Range("C81").Value = -12.837
$ws.Range("D81").Value = -7.929
$ws.Range("B82").Value = 5.375
$ws.Range("B83").Value = 5.351
$ws.Range("C83").Value = -13.97
$ws.Range("E84").Value = 16.734
$ws.Range("D88").Value = -7.972
$ws.Range("B93").Value = 5.508999999999999
$ws.Range("C94").Value = -11.012
$ws.Range("D96").Value = -7.467999999999999
$ws.Range("B97").Value = 5.92
$ws.Range("D99").Value = -7.739999999999999
$ws.Range("C100").Value = -11.876
$ws.Range("E100").Value = 16.883
$ws.Range("E101").Value = 16.683
$ws.Range("C102").Value = -13.55
